# feat: params.py and optimization
#
# The sizing/optimization script (params.py) was re-run, producing updated
# performance numbers. Refresh the computed "Valor" column in the
# Compliance Matrix sheet with the new results, and flag the rows that no
# longer meet their requirement as "NOK" in the "Status" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding a "NN.NN %" style result need to stay plain text - otherwise
# Excel's smart entry would reinterpret them as a percentage number. Force
# those ranges to Text format before writing the new values.
$ws.Range("C10:C15").NumberFormat = "@"
$ws.Range("C19:C20").NumberFormat = "@"
$ws.Range("C30:C31").NumberFormat = "@"
$ws.Range("C35").NumberFormat = "@"

# --- Column C ("Valor") updates ---
$ws.Range("C4").Value = "476052.66 N"

$ws.Range("C10").Value = "10.74 %"
$ws.Range("C11").Value = "7.49 %"
$ws.Range("C12").Value = "7.77 %"
$ws.Range("C13").Value = "13.85 %"
$ws.Range("C14").Value = "25.46 %"
$ws.Range("C15").Value = "11.22 %"

$ws.Range("C18").Value = "26510.29 N (Ref T0 req)"

$ws.Range("C19").Value = "29.85 %"
$ws.Range("C20").Value = "17.92 %"

$ws.Range("C30").Value = "18.07 %"
$ws.Range("C31").Value = "15.17 %"
$ws.Range("C32").Value = "30.10 º"
$ws.Range("C33").Value = "22.45 º"
$ws.Range("C34").Value = "61.50 º"
$ws.Range("C35").Value = "123.24 %"
$ws.Range("C36").Value = "26.68 m"

# --- Column D ("Status") updates: rows that now fail compliance ---
$ws.Range("D30").Value = "NOK"
$ws.Range("D35").Value = "NOK"
